$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "E2"; Value = 1 },
    @{ Cell = "F2"; Value = 0.3333333333333333 },
    @{ Cell = "G2"; Value = 0.1050873333333333 },
    @{ Cell = "H2"; Value = 0.315262 },
    @{ Cell = "I2"; Value = 0.03031434174852429 },
    @{ Cell = "J2"; Value = 0.03031434174852429 },
    @{ Cell = "M2"; Value = 16.28844733333333 },
    @{ Cell = "N2"; Value = 48.865342 },
    @{ Cell = "O2"; Value = 0.2176904746803693 },
    @{ Cell = "P2"; Value = 0.2176904746803693 },
    @{ Cell = "Q2"; Value = 1.711709494400444 },
    @{ Cell = "R2"; Value = 15.405385449604 },
    @{ Cell = "S2"; Value = 0.006599143444859191 },
    @{ Cell = "T2"; Value = 0.006599143444859191 },
    @{ Cell = "E3"; Value = 1 },
    @{ Cell = "F3"; Value = 0.3333333333333333 },
    @{ Cell = "G3"; Value = 0.1050873333333333 },
    @{ Cell = "H3"; Value = 0.315262 },
    @{ Cell = "I3"; Value = 0.03031434174852429 },
    @{ Cell = "J3"; Value = 0.03031434174852429 },
    @{ Cell = "M3"; Value = 27.61090666666666 },
    @{ Cell = "N3"; Value = 82.83271999999999 },
    @{ Cell = "O3"; Value = 0.3690119294748028 },
    @{ Cell = "P3"; Value = 0.3690119294748029 },
    @{ Cell = "Q3"; Value = 2.901556552515555 },
    @{ Cell = "R3"; Value = 26.11400897264 },
    @{ Cell = "S3"; Value = 0.01118635373938152 },
    @{ Cell = "T3"; Value = 0.01118635373938152 },
    @{ Cell = "E4"; Value = 1 },
    @{ Cell = "F4"; Value = 0.3333333333333333 },
    @{ Cell = "G4"; Value = 0.1050873333333333 },
    @{ Cell = "H4"; Value = 0.315262 },
    @{ Cell = "I4"; Value = 0.03031434174852429 },
    @{ Cell = "J4"; Value = 0.03031434174852429 },
    @{ Cell = "M4"; Value = 26.266325 },
    @{ Cell = "N4"; Value = 78.798975 },
    @{ Cell = "O4"; Value = 0.3510419771967738 },
    @{ Cell = "P4"; Value = 0.3510419771967739 },
    @{ Cell = "Q4"; Value = 2.760258050716666 },
    @{ Cell = "R4"; Value = 24.84232245645 },
    @{ Cell = "S4"; Value = 0.01064160646482067 },
    @{ Cell = "T4"; Value = 0.01064160646482067 },
    @{ Cell = "E5"; Value = 1 },
    @{ Cell = "F5"; Value = 0.3333333333333333 },
    @{ Cell = "G5"; Value = 0.1050873333333333 },
    @{ Cell = "H5"; Value = 0.315262 },
    @{ Cell = "I5"; Value = 0.03031434174852429 },
    @{ Cell = "J5"; Value = 0.03031434174852429 },
    @{ Cell = "M5"; Value = 4.658207333333333 },
    @{ Cell = "N5"; Value = 13.974622 },
    @{ Cell = "O5"; Value = 0.06225561864805391 },
    @{ Cell = "P5"; Value = 0.06225561864805392 },
    @{ Cell = "Q5"; Value = 0.4895185867737777 },
    @{ Cell = "R5"; Value = 4.405667280964 },
    @{ Cell = "S5"; Value = 0.001887238099462908 },
    @{ Cell = "T5"; Value = 0.001887238099462908 },
    @{ Cell = "I6"; Value = 0.6354599969768544 },
    @{ Cell = "J6"; Value = 0.6354599969768545 },
    @{ Cell = "M6"; Value = 16.28844733333333 },
    @{ Cell = "N6"; Value = 48.865342 },
    @{ Cell = "O6"; Value = 0.2176904746803693 },
    @{ Cell = "P6"; Value = 0.2176904746803693 },
    @{ Cell = "Q6"; Value = 35.88146228475866 },
    @{ Cell = "R6"; Value = 322.9331605628279 },
    @{ Cell = "S6"; Value = 0.1383335883822775 },
    @{ Cell = "T6"; Value = 0.1383335883822775 },
    @{ Cell = "I7"; Value = 0.6354599969768544 },
    @{ Cell = "J7"; Value = 0.6354599969768545 },
    @{ Cell = "M7"; Value = 27.61090666666666 },
    @{ Cell = "N7"; Value = 82.83271999999999 },
    @{ Cell = "O7"; Value = 0.3690119294748028 },
    @{ Cell = "P7"; Value = 0.3690119294748029 },
    @{ Cell = "Q7"; Value = 60.82345885605332 },
    @{ Cell = "R7"; Value = 547.4111297044799 },
    @{ Cell = "S7"; Value = 0.2344923195884814 },
    @{ Cell = "T7"; Value = 0.2344923195884815 },
    @{ Cell = "I8"; Value = 0.6354599969768544 },
    @{ Cell = "J8"; Value = 0.6354599969768545 },
    @{ Cell = "M8"; Value = 26.266325 },
    @{ Cell = "N8"; Value = 78.798975 },
    @{ Cell = "O8"; Value = 0.3510419771967738 },
    @{ Cell = "P8"; Value = 0.3510419771967739 },
    @{ Cell = "Q8"; Value = 57.86150948334999 },
    @{ Cell = "R8"; Value = 520.7535853501499 },
    @{ Cell = "S8"; Value = 0.2230731337682109 },
    @{ Cell = "T8"; Value = 0.223073133768211 },
    @{ Cell = "I9"; Value = 0.6354599969768544 },
    @{ Cell = "J9"; Value = 0.6354599969768545 },
    @{ Cell = "M9"; Value = 4.658207333333333 },
    @{ Cell = "N9"; Value = 13.974622 },
    @{ Cell = "O9"; Value = 0.06225561864805391 },
    @{ Cell = "P9"; Value = 0.06225561864805392 },
    @{ Cell = "Q9"; Value = 10.26146245403867 },
    @{ Cell = "R9"; Value = 92.35316208634799 },
    @{ Cell = "S9"; Value = 0.03956095523788454 },
    @{ Cell = "T9"; Value = 0.03956095523788455 },
    @{ Cell = "E10"; Value = 3 },
    @{ Cell = "F10"; Value = 1 },
    @{ Cell = "G10"; Value = 0.9666886666666668 },
    @{ Cell = "H10"; Value = 2.900066 },
    @{ Cell = "I10"; Value = 0.2788588279503266 },
    @{ Cell = "J10"; Value = 0.2788588279503266 },
    @{ Cell = "M10"; Value = 16.28844733333333 },
    @{ Cell = "N10"; Value = 48.865342 },
    @{ Cell = "O10"; Value = 0.2176904746803693 },
    @{ Cell = "P10"; Value = 0.2176904746803693 },
    @{ Cell = "Q10"; Value = 15.74585743473022 },
    @{ Cell = "R10"; Value = 141.712716912572 },
    @{ Cell = "S10"; Value = 0.06070491062531805 },
    @{ Cell = "T10"; Value = 0.06070491062531805 },
    @{ Cell = "E11"; Value = 3 },
    @{ Cell = "F11"; Value = 1 },
    @{ Cell = "G11"; Value = 0.9666886666666668 },
    @{ Cell = "H11"; Value = 2.900066 },
    @{ Cell = "I11"; Value = 0.2788588279503266 },
    @{ Cell = "J11"; Value = 0.2788588279503266 },
    @{ Cell = "M11"; Value = 27.61090666666666 },
    @{ Cell = "N11"; Value = 82.83271999999999 },
    @{ Cell = "O11"; Value = 0.3690119294748028 },
    @{ Cell = "P11"; Value = 0.3690119294748029 },
    @{ Cell = "Q11"; Value = 26.69115055105778 },
    @{ Cell = "R11"; Value = 240.22035495952 },
    @{ Cell = "S11"; Value = 0.1029022341530321 },
    @{ Cell = "T11"; Value = 0.1029022341530321 },
    @{ Cell = "E12"; Value = 3 },
    @{ Cell = "F12"; Value = 1 },
    @{ Cell = "G12"; Value = 0.9666886666666668 },
    @{ Cell = "H12"; Value = 2.900066 },
    @{ Cell = "I12"; Value = 0.2788588279503266 },
    @{ Cell = "J12"; Value = 0.2788588279503266 },
    @{ Cell = "M12"; Value = 26.266325 },
    @{ Cell = "N12"; Value = 78.798975 },
    @{ Cell = "O12"; Value = 0.3510419771967738 },
    @{ Cell = "P12"; Value = 0.3510419771967739 },
    @{ Cell = "Q12"; Value = 25.39135869248333 },
    @{ Cell = "R12"; Value = 228.52222823235 },
    @{ Cell = "S12"; Value = 0.09789115432245764 },
    @{ Cell = "T12"; Value = 0.09789115432245765 },
    @{ Cell = "E13"; Value = 3 },
    @{ Cell = "F13"; Value = 1 },
    @{ Cell = "G13"; Value = 0.9666886666666668 },
    @{ Cell = "H13"; Value = 2.900066 },
    @{ Cell = "I13"; Value = 0.2788588279503266 },
    @{ Cell = "J13"; Value = 0.2788588279503266 },
    @{ Cell = "M13"; Value = 4.658207333333333 },
    @{ Cell = "N13"; Value = 13.974622 },
    @{ Cell = "O13"; Value = 0.06225561864805391 },
    @{ Cell = "P13"; Value = 0.06225561864805392 },
    @{ Cell = "Q13"; Value = 4.503036236116889 },
    @{ Cell = "R13"; Value = 40.52732612505201 },
    @{ Cell = "S13"; Value = 0.01736052884951881 },
    @{ Cell = "T13"; Value = 0.01736052884951881 },
    @{ Cell = "G14"; Value = 0.191934 },
    @{ Cell = "H14"; Value = 0.575802 },
    @{ Cell = "I14"; Value = 0.05536683332429468 },
    @{ Cell = "J14"; Value = 0.05536683332429467 },
    @{ Cell = "M14"; Value = 16.28844733333333 },
    @{ Cell = "N14"; Value = 48.865342 },
    @{ Cell = "O14"; Value = 0.2176904746803693 },
    @{ Cell = "P14"; Value = 0.2176904746803693 },
    @{ Cell = "Q14"; Value = 3.126306850476 },
    @{ Cell = "R14"; Value = 28.136761654284 },
    @{ Cell = "S14"; Value = 0.0120528322279146 },
    @{ Cell = "T14"; Value = 0.0120528322279146 },
    @{ Cell = "G15"; Value = 0.191934 },
    @{ Cell = "H15"; Value = 0.575802 },
    @{ Cell = "I15"; Value = 0.05536683332429468 },
    @{ Cell = "J15"; Value = 0.05536683332429467 },
    @{ Cell = "M15"; Value = 27.61090666666666 },
    @{ Cell = "N15"; Value = 82.83271999999999 },
    @{ Cell = "O15"; Value = 0.3690119294748028 },
    @{ Cell = "P15"; Value = 0.3690119294748029 },
    @{ Cell = "Q15"; Value = 5.29947176016 },
    @{ Cell = "R15"; Value = 47.69524584144 },
    @{ Cell = "S15"; Value = 0.02043102199390779 },
    @{ Cell = "T15"; Value = 0.02043102199390779 },
    @{ Cell = "G16"; Value = 0.191934 },
    @{ Cell = "H16"; Value = 0.575802 },
    @{ Cell = "I16"; Value = 0.05536683332429468 },
    @{ Cell = "J16"; Value = 0.05536683332429467 },
    @{ Cell = "M16"; Value = 26.266325 },
    @{ Cell = "N16"; Value = 78.798975 },
    @{ Cell = "O16"; Value = 0.3510419771967738 },
    @{ Cell = "P16"; Value = 0.3510419771967739 },
    @{ Cell = "Q16"; Value = 5.04140082255 },
    @{ Cell = "R16"; Value = 45.37260740295 },
    @{ Cell = "S16"; Value = 0.01943608264128463 },
    @{ Cell = "T16"; Value = 0.01943608264128463 },
    @{ Cell = "G17"; Value = 0.191934 },
    @{ Cell = "H17"; Value = 0.575802 },
    @{ Cell = "I17"; Value = 0.05536683332429468 },
    @{ Cell = "J17"; Value = 0.05536683332429467 },
    @{ Cell = "M17"; Value = 4.658207333333333 },
    @{ Cell = "N17"; Value = 13.974622 },
    @{ Cell = "O17"; Value = 0.06225561864805391 },
    @{ Cell = "P17"; Value = 0.06225561864805392 },
    @{ Cell = "Q17"; Value = 0.8940683663160001 },
    @{ Cell = "R17"; Value = 8.046615296844001 },
    @{ Cell = "S17"; Value = 0.003446896461187653 },
    @{ Cell = "T17"; Value = 0.003446896461187653 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
